$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New package rows to append at the bottom of the table (rows 27-29).
# Cell-write order intentionally follows the order the new strings were
# entered (matches the resulting sharedStrings.xml ordering).

# pkgdown
$ws.Range("A27").Value = "pkgdown"
$ws.Range("B27").Value = "Seus pacotes maravilhosamente documentados."
$ws.Range("C27").Value = "https://pkgdown.r-lib.org/"

# testthat
$ws.Range("A28").Value = "testthat"
$ws.Range("C28").Value = "https://testthat.r-lib.org/"
$ws.Range("B28").Value = "Teste."

# sf
$ws.Range("C29").Value = "https://r-spatial.github.io/sf/"
$ws.Range("B29").Value = "Simplificando dados espaciais."
$ws.Range("A29").Value = "sf"

# Update the view to match the author's final state
$ws.Range("A29").Select()
